$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 66
$ws1.Range("F4").Value = 167
$ws1.Range("F6").Value = 5442
$ws1.Range("F8").Value = 5401
$ws1.Range("F9").Value = 635
$ws1.Range("F11").Value = 1385
$ws1.Range("F12").Value = 14

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 66
$ws4.Range("F4").Value = 167
$ws4.Range("F7").Value = 5442
$ws4.Range("F9").Value = 5401
$ws4.Range("F10").Value = 635
$ws4.Range("F12").Value = 1385
$ws4.Range("F13").Value = 14
